# Updated symbol list on Fri Feb 17 14:25:26 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) quotes for
# the coin rows on Sheet1. Values are stored as plain text in the source
# workbook (e.g. "309.55", "-2.66%"), so each target cell's NumberFormat
# is forced to Text ("@") before the new literal is written - this stops
# Excel's COM layer from auto-converting numeric- or percent-looking
# strings into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}


# Row 2
Set-TextValue $ws "D2" "309.54"
Set-TextValue $ws "E2" "-2.54%"

# Row 3
Set-TextValue $ws "D3" "52.48"
Set-TextValue $ws "E3" "9.62%"

# Row 4
Set-TextValue $ws "D4" "5.110"
Set-TextValue $ws "E4" "-3.27%"

# Row 5
Set-TextValue $ws "D5" "0.07806"
Set-TextValue $ws "E5" "-1.67%"

# Row 6
Set-TextValue $ws "D6" "4.503"
Set-TextValue $ws "E6" "-1.99%"

# Row 7
Set-TextValue $ws "E7" "3.07%"

# Row 8
Set-TextValue $ws "D8" "1.579"
Set-TextValue $ws "E8" "-3.85%"

# Row 9
Set-TextValue $ws "E9" "-3.65%"

# Row 10
Set-TextValue $ws "D10" "0.2001"
Set-TextValue $ws "E10" "3.32%"

# Row 11
Set-TextValue $ws "D11" "0.09593"
Set-TextValue $ws "E11" "2.19%"

# Row 12
Set-TextValue $ws "D12" "0.04719"
Set-TextValue $ws "E12" "1.84%"

# Row 13
Set-TextValue $ws "E13" "0.07%"

# Row 14
Set-TextValue $ws "D14" "0.001261"
Set-TextValue $ws "E14" "-4.62%"

# Row 15
Set-TextValue $ws "D15" "0.005788"
Set-TextValue $ws "E15" "-1.02%"

# Row 16
Set-TextValue $ws "E16" "2,010.46%"

# Row 17
Set-TextValue $ws "D17" "3.332"
Set-TextValue $ws "E17" "0.02%"

# Row 18
Set-TextValue $ws "D18" "2.415"
Set-TextValue $ws "E18" "-0.67%"

# Row 19
Set-TextValue $ws "D19" "0.3445"
Set-TextValue $ws "E19" "-0.62%"

# Row 20
Set-TextValue $ws "D20" "7.994"
Set-TextValue $ws "E20" "-1.25%"

# Row 21
Set-TextValue $ws "D21" "0.1363"
Set-TextValue $ws "E21" "-2.21%"

# Row 22
Set-TextValue $ws "D22" "0.3090"
Set-TextValue $ws "E22" "-0.40%"

# Row 23
Set-TextValue $ws "D23" "0.04175"
Set-TextValue $ws "E23" "0.10%"

# Row 24
Set-TextValue $ws "D24" "0.001260"
Set-TextValue $ws "E24" "-4.64%"

# Row 25
Set-TextValue $ws "D25" "0.003968"
Set-TextValue $ws "E25" "-6.54%"

# Row 26
Set-TextValue $ws "D26" "0.0001348"
Set-TextValue $ws "E26" "-0.41%"

# Row 38
Set-TextValue $ws "D38" "0.02604"
Set-TextValue $ws "E38" "-1.72%"

# Row 39
Set-TextValue $ws "D39" "0.05892"
Set-TextValue $ws "E39" "1.71%"

# Row 40
Set-TextValue $ws "D40" "0.01126"
Set-TextValue $ws "E40" "4.40%"

# Row 41
Set-TextValue $ws "D41" "0.007898"
Set-TextValue $ws "E41" "-1.45%"

# Row 42
Set-TextValue $ws "D42" "0.1426"
Set-TextValue $ws "E42" "-0.49%"

# Row 43
Set-TextValue $ws "D43" "0.008226"
Set-TextValue $ws "E43" "6.88%"

# Row 44
Set-TextValue $ws "E44" "-0.64%"

# Row 45
Set-TextValue $ws "D45" "0.3123"
Set-TextValue $ws "E45" "-1.32%"

# Row 46
Set-TextValue $ws "D46" "0.00007274"
Set-TextValue $ws "E46" "4.93%"

# Row 47
Set-TextValue $ws "D47" "0.00000000749"
Set-TextValue $ws "E47" "-0.36%"

# Row 48
Set-TextValue $ws "D48" "0.05667"
Set-TextValue $ws "E48" "3.35%"

# Row 49
Set-TextValue $ws "D49" "0.002617"

# Row 50
Set-TextValue $ws "D50" "0.00002097"
Set-TextValue $ws "E50" "-0.36%"

# Row 51
Set-TextValue $ws "D51" "0.0001998"
Set-TextValue $ws "E51" "-0.36%"
